# =====================================================================
# Scheduled market-data refresh for Halicarnassus_Profits workbook
# Refreshes cached Universalis price snapshots and recalculated leve
# profit figures (columns H:N) across all job sheets.
# =====================================================================

$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 155.27777
$ws.Range("I9").Value = 228.2
$ws.Range("J9").Value = 64.125
$ws.Range("K9").Value = 228.2
$ws.Range("L9").Value = 64.125
$ws.Range("M9").Value = -59.19999999999999
$ws.Range("N9").Value = -402.125
$ws.Range("H12").Value = 887.5
$ws.Range("H17").Value = 1878.5
$ws.Range("J17").Value = 2223.125
$ws.Range("L17").Value = 6669.375
$ws.Range("N17").Value = -7005.375
$ws.Range("H21").Value = 10700
$ws.Range("I21").Value = 10700
$ws.Range("K21").Value = 10700
$ws.Range("M21").Value = -10232
$ws.Range("H23").Value = 10700
$ws.Range("I23").Value = 10700
$ws.Range("K23").Value = 10700
$ws.Range("M23").Value = -10466
$ws.Range("H29").Value = 3393.3333
$ws.Range("J29").Value = 4475
$ws.Range("L29").Value = 13425
$ws.Range("N29").Value = -13987
$ws.Range("H38").Value = 1882.1
$ws.Range("I38").Value = 1423.5555
$ws.Range("K38").Value = 4270.666499999999
$ws.Range("M38").Value = -3898.666499999999
$ws.Range("H51").Value = 2856.25
$ws.Range("I51").Value = 2000
$ws.Range("K51").Value = 2000
$ws.Range("M51").Value = -1516
$ws.Range("H58").Value = 2170.3076
$ws.Range("I58").Value = 604.6667
$ws.Range("J58").Value = 2640
$ws.Range("K58").Value = 1814.0001
$ws.Range("L58").Value = 7920
$ws.Range("M58").Value = -1664.0001
$ws.Range("N58").Value = -8220
$ws.Range("H88").Value = 877.3333
$ws.Range("I88").Value = 1174.75
$ws.Range("J88").Value = 639.4
$ws.Range("K88").Value = 1174.75
$ws.Range("L88").Value = 639.4
$ws.Range("M88").Value = -768.75
$ws.Range("N88").Value = -1451.4
$ws.Range("H91").Value = 877.3333
$ws.Range("I91").Value = 1174.75
$ws.Range("J91").Value = 639.4
$ws.Range("K91").Value = 1174.75
$ws.Range("L91").Value = 639.4
$ws.Range("M91").Value = 229.25
$ws.Range("N91").Value = -3447.4
$ws.Range("H113").Value = 5620.75
$ws.Range("I113").Value = 2991.5
$ws.Range("J113").Value = 8250
$ws.Range("K113").Value = 2991.5
$ws.Range("L113").Value = 8250
$ws.Range("M113").Value = 262.5
$ws.Range("N113").Value = -14758
$ws.Range("H135").Value = 925.8889
$ws.Range("J135").Value = 484.5
$ws.Range("L135").Value = 4360.5
$ws.Range("N135").Value = -9430.5
$ws.Range("H137").Value = 3330.625
$ws.Range("I137").Value = 2184.5715
$ws.Range("J137").Value = 4222
$ws.Range("K137").Value = 6553.7145
$ws.Range("L137").Value = 12666
$ws.Range("M137").Value = -4003.7145
$ws.Range("N137").Value = -17766
$ws.Range("H138").Value = 2360.2307
$ws.Range("I138").Value = 2390.25
$ws.Range("J138").Value = 2000
$ws.Range("K138").Value = 7170.75
$ws.Range("L138").Value = 6000
$ws.Range("M138").Value = -2030.75
$ws.Range("N138").Value = -16280

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2782.7693
$ws.Range("I45").Value = 2243
$ws.Range("K45").Value = 2243
$ws.Range("M45").Value = -1866
$ws.Range("H61").Value = 4191.7856
$ws.Range("I61").Value = 2711
$ws.Range("K61").Value = 2711
$ws.Range("M61").Value = -2499
$ws.Range("H74").Value = 2075.6206
$ws.Range("I74").Value = 2122.4092
$ws.Range("J74").Value = 1928.5714
$ws.Range("K74").Value = 2122.4092
$ws.Range("L74").Value = 1928.5714
$ws.Range("M74").Value = -1248.4092
$ws.Range("N74").Value = -3676.5714
$ws.Range("H77").Value = 2075.6206
$ws.Range("I77").Value = 2122.4092
$ws.Range("J77").Value = 1928.5714
$ws.Range("K77").Value = 10612.046
$ws.Range("L77").Value = 9642.857
$ws.Range("M77").Value = -6244.046
$ws.Range("N77").Value = -18378.857
$ws.Range("H88").Value = 1325.8125
$ws.Range("J88").Value = 985.7143
$ws.Range("L88").Value = 985.7143
$ws.Range("N88").Value = -1797.7143
$ws.Range("H91").Value = 1325.8125
$ws.Range("J91").Value = 985.7143
$ws.Range("L91").Value = 985.7143
$ws.Range("N91").Value = -3793.7143
$ws.Range("H97").Value = 739.5789
$ws.Range("I97").Value = 717.7646999999999
$ws.Range("J97").Value = 925
$ws.Range("K97").Value = 717.7646999999999
$ws.Range("L97").Value = 925
$ws.Range("M97").Value = -221.7646999999999
$ws.Range("N97").Value = -1917
$ws.Range("H132").Value = 3997.2144
$ws.Range("I132").Value = 3997.2144
$ws.Range("K132").Value = 11991.6432
$ws.Range("M132").Value = -9461.643199999999
$ws.Range("H136").Value = 4191.7856
$ws.Range("I136").Value = 2711
$ws.Range("K136").Value = 8133
$ws.Range("M136").Value = -5583

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 896.4545000000001
$ws.Range("I80").Value = 653.1667
$ws.Range("J80").Value = 1188.4
$ws.Range("K80").Value = 653.1667
$ws.Range("L80").Value = 1188.4
$ws.Range("M80").Value = 344.8333
$ws.Range("N80").Value = -3184.4
$ws.Range("H83").Value = 896.4545000000001
$ws.Range("I83").Value = 653.1667
$ws.Range("J83").Value = 1188.4
$ws.Range("K83").Value = 3265.8335
$ws.Range("L83").Value = 5942
$ws.Range("M83").Value = 1726.1665
$ws.Range("N83").Value = -15926
$ws.Range("H94").Value = 1495.45
$ws.Range("I94").Value = 1511.6666
$ws.Range("J94").Value = 1349.5
$ws.Range("K94").Value = 1511.6666
$ws.Range("L94").Value = 1349.5
$ws.Range("M94").Value = -1060.6666
$ws.Range("N94").Value = -2251.5
$ws.Range("H99").Value = 1720.2
$ws.Range("I99").Value = 1696.6666
$ws.Range("J99").Value = 1755.5
$ws.Range("K99").Value = 1696.6666
$ws.Range("L99").Value = 1755.5
$ws.Range("M99").Value = -198.6666
$ws.Range("N99").Value = -4751.5

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6661.84
$ws.Range("I31").Value = 1097
$ws.Range("J31").Value = 8825.944
$ws.Range("K31").Value = 1097
$ws.Range("L31").Value = 8825.944
$ws.Range("M31").Value = -802
$ws.Range("N31").Value = -9415.944
$ws.Range("H34").Value = 6661.84
$ws.Range("I34").Value = 1097
$ws.Range("J34").Value = 8825.944
$ws.Range("K34").Value = 1097
$ws.Range("L34").Value = 8825.944
$ws.Range("M34").Value = -895
$ws.Range("N34").Value = -9229.944
$ws.Range("H94").Value = 4240.6924
$ws.Range("I94").Value = 2434.625
$ws.Range("J94").Value = 7130.4
$ws.Range("K94").Value = 2434.625
$ws.Range("L94").Value = 7130.4
$ws.Range("M94").Value = -1983.625
$ws.Range("N94").Value = -8032.4
$ws.Range("H132").Value = 3950.3076
$ws.Range("I132").Value = 3635.4
$ws.Range("K132").Value = 10906.2
$ws.Range("M132").Value = -8376.200000000001
$ws.Range("H134").Value = 801.9286
$ws.Range("I134").Value = 801.9286
$ws.Range("K134").Value = 2405.7858
$ws.Range("M134").Value = 129.2142000000003

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 802.1429000000001
$ws.Range("I5").Value = 681.1429000000001
$ws.Range("J5").Value = 923.1429000000001
$ws.Range("K5").Value = 2043.4287
$ws.Range("L5").Value = 2769.4287
$ws.Range("M5").Value = -1931.4287
$ws.Range("N5").Value = -2993.4287
$ws.Range("H37").Value = 109539
$ws.Range("J37").Value = 109539
$ws.Range("L37").Value = 328617
$ws.Range("N37").Value = -328841
$ws.Range("H68").Value = 2007
$ws.Range("J68").Value = 2253.5715
$ws.Range("L68").Value = 6760.7145
$ws.Range("N68").Value = -8382.7145
$ws.Range("H71").Value = 2007
$ws.Range("J71").Value = 2253.5715
$ws.Range("L71").Value = 20282.1435
$ws.Range("N71").Value = -28394.1435
$ws.Range("H132").Value = 1968.5625
$ws.Range("I132").Value = 1889.7
$ws.Range("J132").Value = 2100
$ws.Range("K132").Value = 17007.3
$ws.Range("L132").Value = 18900
$ws.Range("M132").Value = -14477.3
$ws.Range("N132").Value = -23960
$ws.Range("H135").Value = 802.1429000000001
$ws.Range("I135").Value = 681.1429000000001
$ws.Range("J135").Value = 923.1429000000001
$ws.Range("K135").Value = 6130.2861
$ws.Range("L135").Value = 8308.286100000001
$ws.Range("M135").Value = -3595.2861
$ws.Range("N135").Value = -13378.2861

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 2078
$ws.Range("I70").Value = 2078
$ws.Range("K70").Value = 2078
$ws.Range("M70").Value = -1808
$ws.Range("H73").Value = 2078
$ws.Range("I73").Value = 2078
$ws.Range("K73").Value = 2078
$ws.Range("M73").Value = -1142
$ws.Range("H132").Value = 4028.7827
$ws.Range("I132").Value = 3058.5293
$ws.Range("J132").Value = 6777.8335
$ws.Range("K132").Value = 9175.5879
$ws.Range("L132").Value = 20333.5005
$ws.Range("M132").Value = -6645.5879
$ws.Range("N132").Value = -25393.5005

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5582.6
$ws.Range("I7").Value = 4053.4443
$ws.Range("K7").Value = 4053.4443
$ws.Range("M7").Value = -3941.4443
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("H100").Value = 9174.875
$ws.Range("I100").Value = 6699.5
$ws.Range("K100").Value = 6699.5
$ws.Range("M100").Value = -6158.5
$ws.Range("H122").Value = 3995.2
$ws.Range("I122").Value = 3995.2
$ws.Range("K122").Value = 11985.6
$ws.Range("M122").Value = -9535.599999999999
$ws.Range("H126").Value = 5582.6
$ws.Range("I126").Value = 4053.4443
$ws.Range("K126").Value = 12160.3329
$ws.Range("M126").Value = -9690.332900000001
$ws.Range("M16").ClearContents()  # value no longer reported; cell removed from row

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 18362.924
$ws.Range("J41").Value = 19465.555
$ws.Range("L41").Value = 19465.555
$ws.Range("N41").Value = -20245.555
$ws.Range("H126").Value = 4994.5557
$ws.Range("I126").Value = 3291
$ws.Range("K126").Value = 9873
$ws.Range("M126").Value = -7403
